$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append a new row of mail-log data to the Logs sheet (row 3)
$logs.Range("A3").Value = "Zou jij dit even op kunnen pakken?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #1: Zou jij dit even op kunnen pakken?"
$logs.Range("D3").Value = "Planning / Afspraak"
$logs.Range("E3").Value = "Beste klant,`nDank voor uw bericht. Kunt u specifieker zijn over welke kwestie u wenst dat we oppakken? Graag ontvangen we meer details om u verder te kunnen helpen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F3").Value = "2025-08-06 19:27:43"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"
$logs.Range("I3").Value = "Ja"
$logs.Range("J3").Value = "Nee"

# Undo the automatic row-height change caused by the multi-line value
$logs.Rows.Item(3).AutoFit()

# Extend the existing conditional formatting rules so they also cover row 3
foreach ($col in @("D", "G", "H", "I", "J")) {
    $fcs = $logs.Range($col + "2").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($logs.Range($col + "2:" + $col + "3"))
    }
}

# Update the Dashboard summary count for "Planning / Afspraak"
$dashboard.Range("B2").Value = 2
